$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force Word to split a run at a given [start,end) sub-range instead
# of silently re-merging it with neighbouring runs that already share the
# same formatting (this engine normalises same-format adjacent runs whenever
# .Text is assigned). Toggling a character property off/on is enough to make
# the run boundary "stick".
# ---------------------------------------------------------------------------
function Split-Run($rng) {
    $rng.Font.Bold = $false
    $rng.Font.Bold = $true
}

# ---------------------------------------------------------------------------
# 1) Paragraph spacing: add <w:spacing w:line="360" w:lineRule="auto"/> to the
#    "TRIAL CERTIFICATE" paragraph, the "This is to certify that ..."
#    paragraph and the "has passed ..." paragraph.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt -like "TRIAL CERTIFICATE*" -or `
        $txt -like "This is to certify that*" -or `
        $txt -like "*passed the S.S.C*") {
        $para.LineSpacingRule = 5   # wdLineSpaceMultiple
        $para.LineSpacing = 18      # 360/240*12 -> w:line="360" w:lineRule="auto"
    }
}

# ---------------------------------------------------------------------------
# 2) "$name$" -> "_________________________" (stays its own run)
# ---------------------------------------------------------------------------
$p12 = $d.Paragraphs.Item(12)
$full = $p12.Range.Text
$pStart = $p12.Range.Start
$idx = $full.IndexOf("`$name`$")
$nameStart = $pStart + $idx
$nameEnd = $nameStart + 6
$nameRange = $d.Range($nameStart, $nameEnd)
$nameRange.Text = "_________________________"

# Re-split from the preceding "This is to certify that " run (same
# formatting, so it silently merged back into one run above).
$newNameEnd = $nameStart + 25
$splitRange = $d.Range($nameStart, $newNameEnd)
Split-Run $splitRange

# Remove the old "_GoBack" bookmark that used to sit right after "$name$"
# (it moves to the end of the next paragraph, see step 4).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 3) "$roll$" -> "__________" and split the run right in the middle of the
#    "year _______" blank so the tail (with the seat-no blank) becomes its
#    own run.
# ---------------------------------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$full13 = $p13.Range.Text
$p13Start = $p13.Range.Start

$rollIdx = $full13.IndexOf("`$roll`$")
$rollStart = $p13Start + $rollIdx
$rollEnd = $rollStart + 6
$rollRange = $d.Range($rollStart, $rollEnd)
$rollRange.Text = "__________"

$full13b = $p13.Range.Text
$yearIdx = $full13b.IndexOf("year ")
$splitOffset = $yearIdx + 5 + 3
$splitPos = $p13Start + $splitOffset
$tailEnd = $p13.Range.End - 1   # exclude the paragraph mark

$tailRange = $d.Range($splitPos, $tailEnd)
Split-Run $tailRange

# ---------------------------------------------------------------------------
# 4) Re-add the "_GoBack" bookmark at the very end of paragraph 13 (after the
#    new tail run). Adding a bookmark exactly at a paragraph-end position is
#    unreliable, so we temporarily append a placeholder, bookmark just before
#    it, then delete the placeholder again.
# ---------------------------------------------------------------------------
$p13.Range.InsertAfter("ZZZ")
$placeholderStart = $p13.Range.End - 1 - 3
$bmRange = $d.Range($placeholderStart, $placeholderStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholderRange = $d.Range($placeholderStart, $placeholderStart + 3)
$placeholderRange.Text = ""

Write-Output "Done"
